$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "321.01"
Set-TextValue $ws.Range("E2") "8.41%"
Set-TextValue $ws.Range("G2") "3"

# Row 3
Set-TextValue $ws.Range("D3") "47.06"
Set-TextValue $ws.Range("E3") "12.40%"
Set-TextValue $ws.Range("G3") "3"

# Row 4
Set-TextValue $ws.Range("D4") "5.278"
Set-TextValue $ws.Range("E4") "5.56%"
Set-TextValue $ws.Range("G4") "3"

# Row 5
Set-TextValue $ws.Range("D5") "0.08107"
Set-TextValue $ws.Range("E5") "7.92%"
Set-TextValue $ws.Range("G5") "3"

# Row 6
Set-TextValue $ws.Range("D6") "4.564"
Set-TextValue $ws.Range("E6") "3.84%"
Set-TextValue $ws.Range("G6") "3"

# Row 7
Set-TextValue $ws.Range("D7") "1.663"
Set-TextValue $ws.Range("E7") "5.21%"
Set-TextValue $ws.Range("G7") "3"

# Row 8
Set-TextValue $ws.Range("D8") "1.092"
Set-TextValue $ws.Range("E8") "17.83%"
Set-TextValue $ws.Range("G8") "3"

# Row 9
Set-TextValue $ws.Range("D9") "0.1316"
Set-TextValue $ws.Range("E9") "10.88%"
Set-TextValue $ws.Range("G9") "3"

# Row 10
Set-TextValue $ws.Range("D10") "0.1963"
Set-TextValue $ws.Range("E10") "7.50%"
Set-TextValue $ws.Range("G10") "3"

# Row 11
Set-TextValue $ws.Range("D11") "0.09645"
Set-TextValue $ws.Range("E11") "8.02%"
Set-TextValue $ws.Range("G11") "3"

# Row 12
Set-TextValue $ws.Range("D12") "0.04573"
Set-TextValue $ws.Range("E12") "12.64%"
Set-TextValue $ws.Range("G12") "3"

# Row 13
Set-TextValue $ws.Range("E13") "0.09%"
Set-TextValue $ws.Range("G13") "3"

# Row 14
Set-TextValue $ws.Range("D14") "0.001327"
Set-TextValue $ws.Range("E14") "3.68%"
Set-TextValue $ws.Range("G14") "3"

# Row 15
Set-TextValue $ws.Range("D15") "0.005772"
Set-TextValue $ws.Range("E15") "-1.71%"
Set-TextValue $ws.Range("G15") "3"

# Row 16
Set-TextValue $ws.Range("D16") "3.386"
Set-TextValue $ws.Range("E16") "0.90%"
Set-TextValue $ws.Range("G16") "3"

# Row 17
Set-TextValue $ws.Range("E17") "1.50%"
Set-TextValue $ws.Range("G17") "3"

# Row 18
Set-TextValue $ws.Range("D18") "0.3393"
Set-TextValue $ws.Range("E18") "2.37%"
Set-TextValue $ws.Range("G18") "3"

# Row 19
Set-TextValue $ws.Range("D19") "8.184"
Set-TextValue $ws.Range("E19") "3.26%"
Set-TextValue $ws.Range("G19") "3"

# Row 20
Set-TextValue $ws.Range("D20") "0.1387"
Set-TextValue $ws.Range("E20") "-1.76%"
Set-TextValue $ws.Range("G20") "3"

# Row 21
Set-TextValue $ws.Range("D21") "0.3074"
Set-TextValue $ws.Range("E21") "-7.03%"
Set-TextValue $ws.Range("G21") "3"

# Row 22
Set-TextValue $ws.Range("D22") "0.04303"
Set-TextValue $ws.Range("E22") "4.19%"
Set-TextValue $ws.Range("G22") "3"

# Row 23
Set-TextValue $ws.Range("D23") "0.001305"
Set-TextValue $ws.Range("E23") "2.96%"
Set-TextValue $ws.Range("G23") "3"

# Row 24
Set-TextValue $ws.Range("D24") "0.004259"
Set-TextValue $ws.Range("E24") "9.60%"
Set-TextValue $ws.Range("G24") "3"

# Row 25
Set-TextValue $ws.Range("D25") "0.0001347"
Set-TextValue $ws.Range("E25") "9.28%"
Set-TextValue $ws.Range("G25") "3"

# Row 26
Set-TextValue $ws.Range("D26") "0.0003714"
Set-TextValue $ws.Range("E26") "-0.27%"
Set-TextValue $ws.Range("G26") "3"

# Row 27
Set-TextValue $ws.Range("G27") "3"

# Row 28
Set-TextValue $ws.Range("G28") "3"

# Row 29
Set-TextValue $ws.Range("G29") "3"

# Row 30
Set-TextValue $ws.Range("G30") "3"

# Row 31
Set-TextValue $ws.Range("G31") "3"

# Row 32
Set-TextValue $ws.Range("G32") "3"

# Row 33
Set-TextValue $ws.Range("G33") "3"

# Row 34
Set-TextValue $ws.Range("G34") "3"

# Row 35
Set-TextValue $ws.Range("G35") "3"

# Row 36
Set-TextValue $ws.Range("G36") "3"

# Row 37
Set-TextValue $ws.Range("G37") "3"

# Row 38
Set-TextValue $ws.Range("D38") "0.02754"
Set-TextValue $ws.Range("E38") "14.93%"
Set-TextValue $ws.Range("G38") "3"

# Row 39
Set-TextValue $ws.Range("D39") "0.05533"
Set-TextValue $ws.Range("E39") "6.65%"
Set-TextValue $ws.Range("G39") "3"

# Row 40
Set-TextValue $ws.Range("D40") "0.005788"
Set-TextValue $ws.Range("E40") "-8.32%"
Set-TextValue $ws.Range("G40") "3"

# Row 41
Set-TextValue $ws.Range("D41") "0.007760"
Set-TextValue $ws.Range("E41") "-0.92%"
Set-TextValue $ws.Range("G41") "3"

# Row 42
Set-TextValue $ws.Range("D42") "0.1446"
Set-TextValue $ws.Range("E42") "9.27%"
Set-TextValue $ws.Range("G42") "3"

# Row 43
Set-TextValue $ws.Range("D43") "0.007666"
Set-TextValue $ws.Range("E43") "3.31%"
Set-TextValue $ws.Range("G43") "3"

# Row 44
Set-TextValue $ws.Range("D44") "0.008843"
Set-TextValue $ws.Range("E44") "22.59%"
Set-TextValue $ws.Range("G44") "3"

# Row 45
Set-TextValue $ws.Range("D45") "0.3505"
Set-TextValue $ws.Range("E45") "18.83%"
Set-TextValue $ws.Range("G45") "3"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006838"
Set-TextValue $ws.Range("E46") "6.37%"
Set-TextValue $ws.Range("G46") "3"

# Row 47
Set-TextValue $ws.Range("E47") "-0.41%"
Set-TextValue $ws.Range("G47") "3"

# Row 48
Set-TextValue $ws.Range("D48") "0.06059"
Set-TextValue $ws.Range("E48") "73.04%"
Set-TextValue $ws.Range("G48") "3"

# Row 49
Set-TextValue $ws.Range("D49") "0.003992"
Set-TextValue $ws.Range("E49") "-5.15%"
Set-TextValue $ws.Range("G49") "3"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002096"
Set-TextValue $ws.Range("E50") "-0.41%"
Set-TextValue $ws.Range("G50") "3"

# Row 51
Set-TextValue $ws.Range("D51") "0.0001996"
Set-TextValue $ws.Range("E51") "-0.41%"
Set-TextValue $ws.Range("G51") "3"
